$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update sheet name to reflect new "through" date
$ws.Name = "Through 2021-09-17"

# Update the row label for September to reflect the new "through" date
$ws.Range("A10").Value = "September (through 09-17)"

# Update September (row 10) monthly carjacking counts
$ws.Range("C10").Value = 30
$ws.Range("E10").Value = 32
$ws.Range("F10").Value = 39
$ws.Range("G10").Value = 63
$ws.Range("H10").Value = 87

# Update Total (row 11) counts
$ws.Range("C11").Value = 411
$ws.Range("E11").Value = 522
$ws.Range("F11").Value = 388
$ws.Range("G11").Value = 847
$ws.Range("H11").Value = 1157
